# Adds two new "Title and Content" slides at the end of the deck
# (matching ppt/slides/slide5.xml + slide6.xml / sldId 260 + 261 in the
# target OOXML) and fills in their placeholder text.

$p = $ppt.ActivePresentation

# --- Slide 5 (sldId 260): "PATH GENERATION" ------------------------------
$s5 = $p.Slides.Add($p.Slides.Count + 1, 2)   # ppLayoutText -> Title and Content

$s5Title = $s5.Shapes.Item(1).TextFrame.TextRange
$s5Title.Text = "PATH GENERATION"
$s5Title.LanguageID = "nb-NO"

# Content placeholder is left empty on this slide.

# --- Slide 6 (sldId 261): path-generation debugging notes ----------------
$s6 = $p.Slides.Add($p.Slides.Count + 1, 2)   # ppLayoutText -> Title and Content

# Title placeholder stays empty on this slide (matches target OOXML).

# Each entry is the paragraph text plus its outline IndentLevel (PowerPoint's
# IndentLevel is 1-based: 1 == top level / no <a:pPr>, 2 == <a:pPr lvl="1"/>).
$s6Lines = @(
    @{Text = "This path is generated:"; Level = 1},
    @{Text = "Kristiansand,Skien,Rail,1,Skien,Oslo,Rail,1,Oslo,Hamar,Rail,1,Hamar,Oslo,Road,1,Oslo,Hamburg,Road,1"; Level = 2},
    @{Text = "While it should never have been generated!"; Level = 2},
    @{Text = ""; Level = 2},
    @{Text = ""; Level = 2},
    @{Text = "TO DO: run the unimodal path generation and debug this first!"; Level = 1},
    @{Text = "Then check two-modal path bugs"; Level = 2}
)

$s6Shape = $s6.Shapes.Item(2)
$s6Tf = $s6Shape.TextFrame

# Build up the paragraphs one at a time (instead of one multi-line Text
# assignment) and set LanguageID/IndentLevel on each freshly-inserted
# paragraph individually -- this COM host only reliably stamps
# rPr/pPr changes onto a paragraph that was just (re)anchored.
$s6Tf.TextRange.Text = $s6Lines[0].Text
$s6Tf.TextRange.LanguageID = "nb-NO"

for ($i = 1; $i -lt $s6Lines.Count; $i++) {
    $null = $s6Tf.TextRange.InsertAfter([char]13 + $s6Lines[$i].Text)
    $s6Para = $s6Tf.TextRange.Paragraphs($i + 1)
    $s6Para.LanguageID = "nb-NO"
    $s6Para.IndentLevel = $s6Lines[$i].Level
}
